$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace old Sesamy/CTO record with "zoom" + zero numeric placeholders
$ws.Range("A2").Value = "zoom"
$ws.Range("B2:H2").Value = 0

# Row 3: replace old Sesamy/CFO record with "zoom" + zero numeric placeholders
$ws.Range("A3").Value = "zoom"
$ws.Range("B3:H3").Value = 0

# Row 4: new record for nClouds
$ws.Range("A4").Value = "nClouds"
$ws.Range("B4").Value = "CRO"
$ws.Range("C4").Value = "Scott Jensen"
$ws.Range("D4").Value = "scottjensen@nclouds.com"
$ws.Range("E4").Value = "AWS Terraform Kubernetes Ansible Puppet Docker Python Jenkins"
$ws.Range("F4").Value = "5 years"
$ws.Range("G4").Value = "`$124K - `$157K"
$ws.Range("H4").Value = "https://www.indeed.com/rc/clk?jk=c45674b9f68a5d3c&fccid=c76149658a7e6a8d&vjs=3"
